$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as literal text
# (e.g. "318.43", "1.000") in the source data. Excel normally auto-
# converts such strings to numbers on assignment, so for any updated
# Price cell whose new value would parse as a number we first mark the
# cell as Text (@) to keep it stored as a string, matching the source.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.339.04"
$ws.Range("E2").Value = "  +5.35%  "
$ws.Range("D3").Value = "1.810.64"
$ws.Range("E3").Value = "  +4.57%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "318.43"
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "0.5706"
$ws.Range("E7").Value = "  +14.44%  "
$ws.Range("D8").Value = "0.3885"
$ws.Range("E8").Value = "  +10.85%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07600"
$ws.Range("E9").Value = "  +4.82%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "42.81"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  +7.86%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  +6.17%  "
$ws.Range("E14").Value = "  +6.45%  "
$ws.Range("D15").Value = "1.810.51"
$ws.Range("E15").Value = "  +5.18%  "
$ws.Range("D16").Value = "7.274"
$ws.Range("E16").Value = "  +6.70%  "
$ws.Range("D17").Value = "91.96"
$ws.Range("E17").Value = "  +5.91%  "
$ws.Range("D18").Value = "0.00001074"
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").Value = "0.06476"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("E22").Value = "  +4.96%  "
$ws.Range("D23").Value = "28.357.77"
$ws.Range("E23").Value = "  +5.18%  "
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").Value = "2.143"
$ws.Range("E25").Value = "  +4.41%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.463"
$ws.Range("E26").Value = "  +18.20%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "158.16"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "20.80"
$ws.Range("E28").Value = "  +4.27%  "
$ws.Range("D29").Value = "2.020.08"
$ws.Range("E29").Value = "  +5.09%  "
$ws.Range("D30").Value = "124.06"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").Value = "1.160"
$ws.Range("E31").Value = "  +10.81%  "
$ws.Range("E32").Value = "  +13.89%  "
$ws.Range("D33").Value = "5.792"
$ws.Range("E33").Value = "  +7.35%  "
$ws.Range("D34").Value = "3.630"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").Value = "0.2217"
$ws.Range("E35").Value = "  +11.56%  "
$ws.Range("D36").Value = "8.974"
$ws.Range("E36").Value = "  +20.61%  "
$ws.Range("D37").Value = "0.02319"
$ws.Range("E37").Value = "  +6.12%  "
$ws.Range("D38").Value = "11.68"
$ws.Range("E38").Value = "  +6.18%  "
$ws.Range("D39").Value = "0.06134"
$ws.Range("E39").Value = "  +3.78%  "
$ws.Range("D40").Value = "0.6403"
$ws.Range("E40").Value = "  +6.51%  "
$ws.Range("D41").Value = "5.040"
$ws.Range("E41").Value = "  +6.12%  "
$ws.Range("E42").Value = "  +4.73%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("D44").Value = "1.381"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("E45").Value = "  +4.93%  "
$ws.Range("D46").Value = "0.6009"
$ws.Range("E46").Value = "  +6.99%  "
$ws.Range("D47").Value = "3.703"
$ws.Range("E47").Value = "  +3.66%  "
$ws.Range("D48").Value = "122.63"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("D49").Value = "1.951"
$ws.Range("E49").Value = "  +5.96%  "
$ws.Range("D50").Value = "1.149"
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("D51").Value = "0.06875"
$ws.Range("E51").Value = "  +3.44%  "
